$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Week-10 stand-up table: fill in the previously empty MIERCOLES / JUEVES / VIERNES
# cells for Paula's block (rows 10-12) with the notes from this week's meeting.
# Values are written in the same order the author's session introduced them so the
# shared-string table grows in the same sequence.
$ws.Range("F11").Value = "Reunión con el profesor, tomamos en cuenta las apreciaciones."
$ws.Range("G11").Value = "Reunión con Guille y Santi para corregir algunos errores"
$ws.Range("G10").Value = "Reunión con el profesor, asignación de tareas"
$ws.Range("E10").Value = "Asistir a la reunión"
$ws.Range("E11").Value = "Correcciones a mi parte de las tablas de casos de uso"
$ws.Range("F10").Value = "Correcciones a tablas y diagrama."
$ws.Range("E12").Value = "Ninguna"
$ws.Range("F12").Value = "Ninguna"
$ws.Range("G12").Value = "Ninguna"

# Move the active selection to where the author left off editing.
$ws.Range("G24").Select()
